$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1420.8206
$ws.Range("I15").Value = 1420.8206
$ws.Range("K15").Value = 4262.4618
$ws.Range("M15").Value = -4093.4618

# Row 33
$ws.Range("H33").Value = 159.6
$ws.Range("I33").Value = 124.875
$ws.Range("K33").Value = 124.875
$ws.Range("M33").Value = 104.125

# Row 51
$ws.Range("H51").Value = 3450
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3450
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3450
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4418

# Row 53
$ws.Range("H53").Value = 258.2
$ws.Range("I53").Value = 307.5
$ws.Range("J53").Value = 225.33333
$ws.Range("K53").Value = 307.5
$ws.Range("L53").Value = 225.33333
$ws.Range("M53").Value = 329.5
$ws.Range("N53").Value = -1499.33333

# Row 58
$ws.Range("H58").Value = 538.3333
$ws.Range("I58").Value = 538.3333
$ws.Range("K58").Value = 1614.9999
$ws.Range("M58").Value = -1464.9999

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# Row 80
$ws.Range("H80").Value = 202.84
$ws.Range("I80").Value = 269.8
$ws.Range("J80").Value = 158.2
$ws.Range("K80").Value = 809.4000000000001
$ws.Range("L80").Value = 474.6
$ws.Range("M80").Value = 188.5999999999999
$ws.Range("N80").Value = -2470.6

# Row 83
$ws.Range("H83").Value = 202.84
$ws.Range("I83").Value = 269.8
$ws.Range("J83").Value = 158.2
$ws.Range("K83").Value = 2428.2
$ws.Range("L83").Value = 1423.8
$ws.Range("M83").Value = 2563.8
$ws.Range("N83").Value = -11407.8

# Row 135
$ws.Range("H135").Value = 1389.6111
$ws.Range("J135").Value = 1520.625
$ws.Range("L135").Value = 13685.625
$ws.Range("N135").Value = -18755.625

# Row 137
$ws.Range("H137").Value = 2042.6471
$ws.Range("I137").Value = 1909.0714
$ws.Range("K137").Value = 5727.2142
$ws.Range("M137").Value = -3177.2142

# Row 138
$ws.Range("H138").Value = 4470.9287
$ws.Range("I138").Value = 1127.1111
$ws.Range("J138").Value = 6054.8423
$ws.Range("K138").Value = 3381.3333
$ws.Range("L138").Value = 18164.5269
$ws.Range("M138").Value = 1758.6667
$ws.Range("N138").Value = -28444.5269

$ws = $wb.Worksheets.Item("ARM")
# Row 35
$ws.Range("H35").Value = 1850.25
$ws.Range("I35").Value = 1850.25
$ws.Range("K35").Value = 1850.25
$ws.Range("M35").Value = -1444.25

# Row 61
$ws.Range("H61").Value = 2047.1428
$ws.Range("I61").Value = 1847.9412
$ws.Range("K61").Value = 1847.9412
$ws.Range("M61").Value = -1635.9412

# Row 102
$ws.Range("H102").Value = 862.5
$ws.Range("J102").Value = 500
$ws.Range("L102").Value = 500
$ws.Range("N102").Value = -3744

# Row 132
$ws.Range("H132").Value = 2077.625
$ws.Range("J132").Value = 2088.4
$ws.Range("L132").Value = 6265.200000000001
$ws.Range("N132").Value = -11325.2

# Row 136
$ws.Range("H136").Value = 2047.1428
$ws.Range("I136").Value = 1847.9412
$ws.Range("K136").Value = 5543.8236
$ws.Range("M136").Value = -2993.8236

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2298.92
$ws.Range("I134").Value = 2010.3334
$ws.Range("J134").Value = 3041
$ws.Range("K134").Value = 6031.0002
$ws.Range("L134").Value = 9123
$ws.Range("M134").Value = -3496.0002
$ws.Range("N134").Value = -14193

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 336.14285
$ws.Range("I22").Value = 292.16666
$ws.Range("K22").Value = 292.16666
$ws.Range("M22").Value = 57.83334000000002

# Row 31
$ws.Range("H31").Value = 4186.1333
$ws.Range("I31").Value = 3111.2632
$ws.Range("J31").Value = 6042.727
$ws.Range("K31").Value = 3111.2632
$ws.Range("L31").Value = 6042.727
$ws.Range("M31").Value = -2816.2632
$ws.Range("N31").Value = -6632.727

# Row 34
$ws.Range("H34").Value = 4186.1333
$ws.Range("I34").Value = 3111.2632
$ws.Range("J34").Value = 6042.727
$ws.Range("K34").Value = 3111.2632
$ws.Range("L34").Value = 6042.727
$ws.Range("M34").Value = -2909.2632
$ws.Range("N34").Value = -6446.727

# Row 58
$ws.Range("H58").Value = 2795.2646
$ws.Range("I58").Value = 1269.5
$ws.Range("J58").Value = 4974.9287
$ws.Range("K58").Value = 1269.5
$ws.Range("L58").Value = 4974.9287
$ws.Range("M58").Value = -1066.5
$ws.Range("N58").Value = -5380.9287

# Row 134
$ws.Range("H134").Value = 2489.975
$ws.Range("I134").Value = 2040.4286
$ws.Range("K134").Value = 6121.2858
$ws.Range("M134").Value = -3586.2858

# Row 136
$ws.Range("H136").Value = 2795.2646
$ws.Range("I136").Value = 1269.5
$ws.Range("J136").Value = 4974.9287
$ws.Range("K136").Value = 3808.5
$ws.Range("L136").Value = 14924.7861
$ws.Range("M136").Value = -1258.5
$ws.Range("N136").Value = -20024.7861

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 833.5
$ws.Range("I11").Value = 778.3333
$ws.Range("J11").Value = 999
$ws.Range("K11").Value = 2334.9999
$ws.Range("L11").Value = 2997
$ws.Range("M11").Value = -2194.9999
$ws.Range("N11").Value = -3277

# Row 132
$ws.Range("H132").Value = 10844.77
$ws.Range("I132").Value = 10898.3
$ws.Range("J132").Value = 10666.333
$ws.Range("K132").Value = 98084.7
$ws.Range("L132").Value = 95996.997
$ws.Range("M132").Value = -95554.7
$ws.Range("N132").Value = -101056.997

# Row 140
$ws.Range("H140").Value = 1347.2307
$ws.Range("I140").Value = 1182.5834
$ws.Range("K140").Value = 3547.7502
$ws.Range("M140").Value = 1632.2498

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 297182.84
$ws.Range("J11").Value = 402756
$ws.Range("L11").Value = 402756
$ws.Range("N11").Value = -403034

# Row 54
$ws.Range("H54").Value = 14000
$ws.Range("J54").Value = 14000
$ws.Range("L54").Value = 14000
$ws.Range("N54").Value = -14780

# Row 70
$ws.Range("H70").Value = 7466.5
$ws.Range("I70").Value = 7499
$ws.Range("K70").Value = 7499
$ws.Range("M70").Value = -7229

# Row 73
$ws.Range("H73").Value = 7466.5
$ws.Range("I73").Value = 7499
$ws.Range("K73").Value = 7499
$ws.Range("M73").Value = -6563

# Row 123
$ws.Range("H123").Value = 32054.9
$ws.Range("J123").Value = 32054.9
$ws.Range("L123").Value = 32054.9
$ws.Range("N123").Value = -36954.9

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 132
$ws.Range("H132").Value = 2192.1714
$ws.Range("I132").Value = 1825.4
$ws.Range("J132").Value = 2467.25
$ws.Range("K132").Value = 5476.200000000001
$ws.Range("L132").Value = 7401.75
$ws.Range("M132").Value = -2946.200000000001
$ws.Range("N132").Value = -12461.75

$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 14137507
$ws.Range("J42").Value = 14137507
$ws.Range("L42").Value = 14137507
$ws.Range("N42").Value = -14138633

# Row 46
$ws.Range("H46").Value = 2804.125
$ws.Range("I46").Value = 2047.1177
$ws.Range("K46").Value = 2047.1177
$ws.Range("M46").Value = -1859.1177

# Row 49
$ws.Range("H49").Value = 14137507
$ws.Range("J49").Value = 14137507
$ws.Range("L49").Value = 14137507
$ws.Range("N49").Value = -14137801

# Row 61
$ws.Range("H61").Value = 1999
$ws.Range("I61").Value = 1999
$ws.Range("K61").Value = 1999
$ws.Range("M61").Value = -1797

# Row 113
$ws.Range("H113").Value = 1999
$ws.Range("I113").Value = 1999
$ws.Range("K113").Value = 1999
$ws.Range("M113").Value = 171

# Row 124
$ws.Range("H124").Value = 57499.5
$ws.Range("J124").Value = 57499.5
$ws.Range("L124").Value = 57499.5
$ws.Range("N124").Value = -67319.5

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 999999.75
$ws.Range("I10").Value = 999999.75
$ws.Range("K10").Value = 999999.75
$ws.Range("M10").Value = -999830.75

# Row 122
$ws.Range("H122").Value = 4241.6665
$ws.Range("I122").Value = 4241.6665
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12724.9995
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10274.9995
$ws.Range("N122").ClearContents()
